# Auto-generated script applying numeric cell updates across 8 leve-profit
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Values mirror a
# scheduled pricing-data refresh: H/I/J (market prices), K/L (leve totals),
# M/N (profit) columns. A few rows also gain/lose a cell (e.g. a previously
# blank profit cell now has a value, or a stale one is cleared).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 175805.72
$ws.Range("I15").Value = 175805.72
$ws.Range("K15").Value = 527417.16
$ws.Range("M15").Value = -527248.16
$ws.Range("H28").Value = 695068.5
$ws.Range("I28").Value = 1010391.8
$ws.Range("J28").Value = 1357.2
$ws.Range("K28").Value = 1010391.8
$ws.Range("L28").Value = 1357.2
$ws.Range("M28").Value = -1009906.8
$ws.Range("N28").Value = -2327.2
$ws.Range("H41").Value = 5050741.5
$ws.Range("I41").Value = 8547282
$ws.Range("J41").Value = 183.22223
$ws.Range("K41").Value = 8547282
$ws.Range("L41").Value = 183.22223
$ws.Range("M41").Value = -8546842
$ws.Range("N41").Value = -1063.22223
$ws.Range("H107").Value = 1011065.2
$ws.Range("I107").Value = 5556555.5
$ws.Range("J107").Value = 956.2222
$ws.Range("K107").Value = 5556555.5
$ws.Range("L107").Value = 956.2222
$ws.Range("M107").Value = -5554635.5
$ws.Range("N107").Value = -4796.2222
$ws.Range("H125").Value = 18685580
$ws.Range("I125").Value = 466
$ws.Range("J125").Value = 28028136
$ws.Range("K125").Value = 4194
$ws.Range("L125").Value = 252253224
$ws.Range("M125").Value = -1734
$ws.Range("N125").Value = -252258144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6388
$ws.Range("I74").Value = 1149.8572
$ws.Range("J74").Value = 13721.4
$ws.Range("K74").Value = 1149.8572
$ws.Range("L74").Value = 13721.4
$ws.Range("M74").Value = -275.8571999999999
$ws.Range("N74").Value = -15469.4
$ws.Range("H77").Value = 6388
$ws.Range("I77").Value = 1149.8572
$ws.Range("J77").Value = 13721.4
$ws.Range("K77").Value = 5749.286
$ws.Range("L77").Value = 68607
$ws.Range("M77").Value = -1381.286
$ws.Range("N77").Value = -77343
$ws.Range("H102").Value = 930.5
$ws.Range("I102").Value = 930.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 930.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 691.5
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 3976.3462
$ws.Range("I122").Value = 3386.5625
$ws.Range("J122").Value = 4920
$ws.Range("K122").Value = 10159.6875
$ws.Range("L122").Value = 14760
$ws.Range("M122").Value = -7709.6875
$ws.Range("N122").Value = -19660
$ws.Range("H132").Value = 2262.9124
$ws.Range("I132").Value = 1871.826
$ws.Range("J132").Value = 3898.3635
$ws.Range("K132").Value = 5615.478
$ws.Range("L132").Value = 11695.0905
$ws.Range("M132").Value = -3085.478
$ws.Range("N132").Value = -16755.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1783.3334
$ws.Range("I86").Value = 2433.3333
$ws.Range("J86").Value = 1458.3334
$ws.Range("K86").Value = 2433.3333
$ws.Range("L86").Value = 1458.3334
$ws.Range("M86").Value = -1310.3333
$ws.Range("N86").Value = -3704.3334
$ws.Range("H89").Value = 1783.3334
$ws.Range("I89").Value = 2433.3333
$ws.Range("J89").Value = 1458.3334
$ws.Range("K89").Value = 12166.6665
$ws.Range("L89").Value = 7291.666999999999
$ws.Range("M89").Value = -6550.666499999999
$ws.Range("N89").Value = -18523.667
$ws.Range("H94").Value = 640.10345
$ws.Range("I94").Value = 549.4783
$ws.Range("K94").Value = 549.4783
$ws.Range("M94").Value = -98.47829999999999
$ws.Range("H105").Value = 229952.39
$ws.Range("I105").Value = 2429.4849
$ws.Range("K105").Value = 2429.4849
$ws.Range("M105").Value = -682.4848999999999
$ws.Range("H107").Value = 701.3333
$ws.Range("I107").Value = 341.6
$ws.Range("K107").Value = 341.6
$ws.Range("M107").Value = 1578.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1160
$ws.Range("I16").Value = 1111.4286
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1111.4286
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -824.4286
$ws.Range("N16").Value = -2074
$ws.Range("H113").Value = 1160
$ws.Range("I113").Value = 1111.4286
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1111.4286
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1058.5714
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 1503.8
$ws.Range("I122").Value = 1381
$ws.Range("J122").Value = 1585.6666
$ws.Range("K122").Value = 4143
$ws.Range("L122").Value = 4756.9998
$ws.Range("M122").Value = -1693
$ws.Range("N122").Value = -9656.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4100
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 4875
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 14625
$ws.Range("M69").Value = -2189
$ws.Range("N69").Value = -16247
$ws.Range("H72").Value = 4100
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 4875
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 43875
$ws.Range("M72").Value = -4944
$ws.Range("N72").Value = -51987
$ws.Range("H134").Value = 7060.722
$ws.Range("I134").Value = 3849.4
$ws.Range("J134").Value = 11074.875
$ws.Range("K134").Value = 11548.2
$ws.Range("L134").Value = 33224.625
$ws.Range("M134").Value = -6478.200000000001
$ws.Range("N134").Value = -43364.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 63.714287
$ws.Range("I2").Value = 104
$ws.Range("J2").Value = 47.6
$ws.Range("K2").Value = 104
$ws.Range("L2").Value = 47.6
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = -273.6
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H122").Value = 2437.4614
$ws.Range("I122").Value = 2212.375
$ws.Range("J122").Value = 2797.6
$ws.Range("K122").Value = 6637.125
$ws.Range("L122").Value = 8392.799999999999
$ws.Range("M122").Value = -4187.125
$ws.Range("N122").Value = -13292.8
$ws.Range("H126").Value = 2603.4736
$ws.Range("I126").Value = 2275.4666
$ws.Range("J126").Value = 2817.3914
$ws.Range("K126").Value = 6826.399800000001
$ws.Range("L126").Value = 8452.174199999999
$ws.Range("M126").Value = -4356.399800000001
$ws.Range("N126").Value = -13392.1742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2230.4666
$ws.Range("I61").Value = 1606.4445
$ws.Range("J61").Value = 3166.5
$ws.Range("K61").Value = 1606.4445
$ws.Range("L61").Value = 3166.5
$ws.Range("M61").Value = -1404.4445
$ws.Range("N61").Value = -3570.5
$ws.Range("H68").Value = 2023.75
$ws.Range("I68").Value = 1758
$ws.Range("J68").Value = 2466.6667
$ws.Range("K68").Value = 1758
$ws.Range("L68").Value = 2466.6667
$ws.Range("M68").Value = -1009
$ws.Range("N68").Value = -3964.6667
$ws.Range("H71").Value = 2023.75
$ws.Range("I71").Value = 1758
$ws.Range("J71").Value = 2466.6667
$ws.Range("K71").Value = 8790
$ws.Range("L71").Value = 12333.3335
$ws.Range("M71").Value = -5046
$ws.Range("N71").Value = -19821.3335
$ws.Range("H113").Value = 2230.4666
$ws.Range("I113").Value = 1606.4445
$ws.Range("J113").Value = 3166.5
$ws.Range("K113").Value = 1606.4445
$ws.Range("L113").Value = 3166.5
$ws.Range("M113").Value = 563.5554999999999
$ws.Range("N113").Value = -7506.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H122").Value = 935.7
$ws.Range("I122").Value = 1019.625
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 3058.875
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -608.875
$ws.Range("N122").Value = -6700

Write-Host "Applied 198 cell updates across 8 sheets"
